# Fruta / hortaliza, semanal
#
# Inserts one new weekly price record for "Cebollín" (Terminal La Palmera
# de La Serena) above the existing row 147, pushing the remaining records
# (old rows 147-264) down by one row (to 148-265). The new row carries the
# same constant dimension values (market/region/category/etc.) as every
# other row in this sheet, with fresh date + price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 147; everything below shifts down one row
# (old row 147 -> 148, ..., old row 264 -> 265).
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new record.
$ws.Cells.Item(147, 1).Value = 8
$ws.Cells.Item(147, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(147, 3).Value = "Coquimbo"
$ws.Cells.Item(147, 4).Value = 44907
$ws.Cells.Item(147, 5).Value = 4
$ws.Cells.Item(147, 6).Value = 100112037
$ws.Cells.Item(147, 7).Value = "Cebollín"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 1000
$ws.Cells.Item(147, 11).Value = 1200
$ws.Cells.Item(147, 12).Value = 1400
$ws.Cells.Item(147, 13).Value = 1300
$ws.Cells.Item(147, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(147, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(147, 16).Value = 217
$ws.Cells.Item(147, 17).Value = 6
$ws.Cells.Item(147, 18).Value = "Hortaliza"
